$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect before editing
$ws.Unprotect()

# Update the confidentiality notice date from 2021-03-26 to 2021-03-29
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-10
$ws.Range("D2").Value = 0.1014105218424303
$ws.Range("E2").Value = -0.01797115157247575

$ws.Range("D3").Value = 0.1043287944470134
$ws.Range("E3").Value = -0.006628369421122215

$ws.Range("D4").Value = 0.1174430363269229
$ws.Range("E4").Value = 0.001229407425620721

$ws.Range("D5").Value = 0.1378912091682604
$ws.Range("E5").Value = -0.007036975349766772

$ws.Range("D6").Value = 0.1349707412310079
$ws.Range("E6").Value = 0.005467625899280515

$ws.Range("D7").Value = 0.1460510918286074
$ws.Range("E7").Value = -0.007672882672882575

$ws.Range("D8").Value = 0.1291446660699572
$ws.Range("E8").Value = -0.01082706766917296

$ws.Range("D9").Value = 0.1287599390858005
$ws.Range("E9").Value = -0.006343552806636454

$ws.Range("E10").Value = -0.00593766223191794

# Restore sheet protection (its content/state is unchanged by this edit)
$ws.Protect()

$wb.Save()
